# Textbox response formatting fix
# Update task-order sheet names and their stimulus-file cell values.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO-... ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511686746874921"
$ws1.Range("B2").Value = "go_stims-16511686746601162.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686746711624.csv"
$ws1.Range("B4").Value = "go_stims-16511686746726055.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686746864629.csv"

# --- Sheet 2: NB_TO-... ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511686775216618"
$ws2.Range("B2").Value = "TB-1651168676219393.csv"
$ws2.Range("B3").Value = "ZB-match_5-16511686747080488.csv"
$ws2.Range("B4").Value = "ZB-match_0-165116867493402.csv"
$ws2.Range("B5").Value = "OB-16511686752069519.csv"
$ws2.Range("B6").Value = "ZB-match_2-16511686751080635.csv"
$ws2.Range("B7").Value = "OB-16511686754016075.csv"
$ws2.Range("B8").Value = "TB-165116867750351.csv"
$ws2.Range("B9").Value = "TB-16511686762375462.csv"
$ws2.Range("B10").Value = "OB-1651168675328998.csv"

# --- Sheet 3: RS_TO-... (name change only, cell values unchanged) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16511686775226195"

# --- Sheet 4: TOL_TO-... ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511686775702078"
$ws4.Range("B2").Value = "MM_stims-16511686775372338.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686775247355.csv"
$ws4.Range("B4").Value = "MM_stims-16511686775532942.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686775372338.csv"
$ws4.Range("B6").Value = "MM_stims-1651168677569241.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686775542948.csv"

# --- Sheet 5: vSAT_TO-... ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511686776462946"
$ws5.Range("B2").Value = "SAT_stims-1651168677599605.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686776148975.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511686776302817.csv"
$ws5.Range("B5").Value = "SAT_stims-16511686775722518.csv"
